$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# CED sheet: new "prediction" columns (A & B), refreshed C/D values, new rows
# ---------------------------------------------------------------------------
$ced = $wb.Worksheets.Item("CED")

# Header row: B1/D1 change from "... results" to "... prediction"
$ced.Cells.Item(1,2).Value = "validation prediction"
$ced.Cells.Item(1,4).Value = "testing prediction"

# Row 2
$ced.Cells.Item(2,1).Value = 75.3343091929
$ced.Cells.Item(2,2).Value = 85.4673767089844
$ced.Cells.Item(2,3).Value = 110.17365474153
$ced.Cells.Item(2,4).Value = 99.5393142700195
$ced.Cells.Item(2,5).Value = 0.52
$ced.Cells.Item(2,6).Value = 39.86

# Row 3
$ced.Cells.Item(3,1).Value = 38.1188930719
$ced.Cells.Item(3,2).Value = 59.1052856445313
$ced.Cells.Item(3,3).Value = 52.258299428173
$ced.Cells.Item(3,4).Value = 30.3464260101318

# Row 4
$ced.Cells.Item(4,1).Value = 99.9226218037
$ced.Cells.Item(4,2).Value = 86.4268112182617
$ced.Cells.Item(4,3).Value = 91.71323266061
$ced.Cells.Item(4,4).Value = 47.0205917358399

# Row 5
$ced.Cells.Item(5,1).Value = 73.3913166371
$ced.Cells.Item(5,2).Value = 55.6705207824707
$ced.Cells.Item(5,3).Value = 232.88550969278
$ced.Cells.Item(5,4).Value = 131.240203857422

# Row 6
$ced.Cells.Item(6,1).Value = 86.6489246683
$ced.Cells.Item(6,2).Value = 82.2365036010742
$ced.Cells.Item(6,3).Value = 288.14541835177
$ced.Cells.Item(6,4).Value = 180.660690307617

# Row 7
$ced.Cells.Item(7,1).Value = 97.6407610572
$ced.Cells.Item(7,2).Value = 85.548942565918
$ced.Cells.Item(7,3).Value = 73.479614466793
$ced.Cells.Item(7,4).Value = 101.070121765137

# Row 8
$ced.Cells.Item(8,1).Value = 101.3874286238
$ced.Cells.Item(8,2).Value = 71.151237487793
$ced.Cells.Item(8,3).Value = 96.36471013916
$ced.Cells.Item(8,4).Value = 151.055953979492

# Row 9
$ced.Cells.Item(9,1).Value = 111.9314793025
$ced.Cells.Item(9,2).Value = 54.3901138305664
$ced.Cells.Item(9,3).Value = 50.591904210414
$ced.Cells.Item(9,4).Value = 89.4068069458008

# Row 10
$ced.Cells.Item(10,1).Value = 87.6180861663
$ced.Cells.Item(10,2).Value = 63.2286949157715
$ced.Cells.Item(10,3).Value = 83.76170589193
$ced.Cells.Item(10,4).Value = 118.966529846191

# Row 11
$ced.Cells.Item(11,1).Value = 139.8541034812
$ced.Cells.Item(11,2).Value = 111.470520019531
$ced.Cells.Item(11,3).Value = 76.722590928303
$ced.Cells.Item(11,4).Value = 59.7607841491699

# New rows 12-17 (only A & B columns populated)
$ced.Cells.Item(12,1).Value = 19.9327807577
$ced.Cells.Item(12,2).Value = 53.6841278076172

$ced.Cells.Item(13,1).Value = 22.7321777434
$ced.Cells.Item(13,2).Value = 56.8656845092773

$ced.Cells.Item(14,1).Value = 50.0834627166
$ced.Cells.Item(14,2).Value = 49.0335197448731

$ced.Cells.Item(15,1).Value = 195.8674698347
$ced.Cells.Item(15,2).Value = 129.831283569336

$ced.Cells.Item(16,1).Value = 66.1856214908
$ced.Cells.Item(16,2).Value = 82.0281448364258

$ced.Cells.Item(17,1).Value = 87.3141084414
$ced.Cells.Item(17,2).Value = 82.417610168457

# ---------------------------------------------------------------------------
# Selection / view state: the author had range D41:D42 highlighted while
# reviewing each sheet (on top of the sheet's previous cursor position).
# ---------------------------------------------------------------------------
$acid = $wb.Worksheets.Item("acidification")
$null = $acid.Range("B2").Select()
$null = $acid.Range("D41:D42").Select()

$gwp = $wb.Worksheets.Item("GWP")
$null = $gwp.Range("A1").Select()
$null = $gwp.Range("D41:D42").Select()

$ei99 = $wb.Worksheets.Item("EI99")
$null = $ei99.Range("A1").Select()
$null = $ei99.Range("D41:D42").Select()

$hh = $wb.Worksheets.Item("humanhealth")
$null = $hh.Range("A1").Select()
$null = $hh.Range("D41:D42").Select()

$eq = $wb.Worksheets.Item("ecosystemquality")
$null = $eq.Range("A1").Select()
$null = $eq.Range("D41:D42").Select()

# CED keeps tab focus/selection on D41:D42 (matches the diff's tabSelected sheet)
$null = $ced.Select()
$null = $ced.Range("D41:D42").Select()
